# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) the handback report now
# carries the "target" and "handback" file columns (E/F) in addition to the
# original handoff columns, the status text moves from "Ready for handoff"
# to "Handed back: in sync with en-US", the handback datetime is recorded
# in column G and the handoff reason ("Include") shifts into column H.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # OLE (BGR) form of RGB 6495ED - matches the workbook's HyperLink style

# Overview sheet: status text is shared with the per-language sheets, update
# the two data rows (a.md / b.md) for both the zh-cn and de-de columns.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = $newStatus
$overview.Cells.Item(2, 3).Value = $newStatus
$overview.Cells.Item(3, 2).Value = $newStatus
$overview.Cells.Item(3, 3).Value = $newStatus

# Per-language sheet details: target-file URL, handback-file display name +
# URL, and the new "Latest Handback DateTime" value.
$sheetsInfo = @(
    @{
        Name = "zh-cn";
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/a.md";
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf";
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f14ae3be5ed4e70f0f19c6ac1bfb67c49a3a9ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf";
        HandbackDateTime = "2016-02-23 08:49:27"
    },
    @{
        Name = "de-de";
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/a.md";
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf";
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a370d7c138e72ff3091517048e980fa54bb17dbf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf";
        HandbackDateTime = "2016-02-23 08:49:47"
    }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    foreach ($r in 2, 3) {
        # B: Status -> handed back
        $ws.Cells.Item($r, 2).Value = $newStatus

        # E: Latest Target File (mirrors the "a.md" source-file hyperlink)
        $eCell = $ws.Cells.Item($r, 5)
        $eCell.Value = "a.md"
        $ws.Hyperlinks.Add($eCell, $info.MdUrl, "", "", "a.md")
        $eCell.Font.Underline = 2
        $eCell.Font.Color = $hyperlinkColor

        # F: Latest Handback File (the translated .xlf that came back)
        $fCell = $ws.Cells.Item($r, 6)
        $fCell.Value = $info.XlfDisplay
        $ws.Hyperlinks.Add($fCell, $info.XlfUrl, "", "", $info.XlfDisplay)
        $fCell.Font.Underline = 2
        $fCell.Font.Color = $hyperlinkColor

        # G: Latest Handback DateTime
        $ws.Cells.Item($r, 7).Value = $info.HandbackDateTime

        # H: Handoff Reason (now correctly populated instead of the stale date)
        $ws.Cells.Item($r, 8).Value = "Include"
    }
}

Write-Host "Handback report generated for zh-cn and de-de"
